# Tambah fitur cetak resi, otomatis status pengiriman, dan PDF download
# Adds shipping-status columns (Status_Pengiriman, No_Resi, Ekspedisi) to the
# data_kendaraan sheet, extending the table from A1:L2 to A1:O2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1) ---
$ws.Range("M1").Value = "Status_Pengiriman"
$ws.Range("N1").Value = "No_Resi"
$ws.Range("O1").Value = "Ekspedisi"

# --- New data cells (row 2) ---
$ws.Range("M2").Value = "Diproses"
$ws.Range("N2").Value = "RESI426356"
$ws.Range("O2").Value = "J&T"

# Match the existing header formatting (bold, centered, thin border) used by
# the other header cells (e.g. L1 "Status") by copying its format onto the
# newly added header cells.
$ws.Range("L1").Copy() | Out-Null
$ws.Range("M1:O1").PasteSpecial(-4122) | Out-Null
